# Adding in demo for fitting two condition pCa curves
# The "force_error" column (F) is removed entirely (header + values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the force_error column (F2:F13 values, and F1 header text)
$ws.Range("F2:F13").ClearContents()
$ws.Range("F1").Value = $null

# Update the selection to match the final state (F1:F1048576 selected, active cell F1)
$ws.Range("F1:F1048576").Select()
